$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$badChar = [char]194
$plusMinus = [char]177
$bad = "$badChar$plusMinus"
$good = "$plusMinus"

for ($r = 2; $r -le 17; $r++) {
    foreach ($col in @("B", "C", "D")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value2
        if ($val -ne $null -and $val.GetType().Name -eq "String" -and $val.Contains($bad)) {
            $cell.Value = $val.Replace($bad, $good)
        }
    }
}
